# "started preparing for the report"
# Adds a beta*gamma column (D) to the existing velocity table (rows 36-40)
# and a new small calculation block (rows 42-48) for dE/dx|min.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B is widened to fit the new "dE/dx|min (MeVg^{-1}cm^2" label.
$ws.Columns.Item(2).ColumnWidth = 25.8571428571429

# --- Existing table (rows 36-40): add column D = beta*gamma ---
$ws.Range("D36").Value = "βγ"
$ws.Range("D37").Formula = '=C37*B37/$A$14'
$ws.Range("D38").Formula = '=C38*B38/$A$14'
$ws.Range("D39").Formula = '=C39*B39/$A$14'
$ws.Range("D40").Formula = '=C40*B40/$A$14'

# --- New block (rows 42-48): dE/dx|min = D(2) * distance * density ---
$ws.Range("B42").Value = "D(2)"

$ws.Range("B43").Value = "dE/dx|min (MeVg^{-1}cm^2"
$ws.Range("B45").Value = "density (gcm^{-3}) "
$ws.Range("B47").Value = "距離 (cm) "
$ws.Range("D43").Value = "(MeV)"

$ws.Range("B44").Value = 1.519
$ws.Range("B46").Value = 1.396
$ws.Range("B48").Value = 2.3
$ws.Range("C43").Formula = '=B44*B48*B46'

[void]$ws.Range("D44").Select()
